$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 1.02
    "C2" = 1.030720824714329
    "D2" = 1.039359687964525
    "E2" = 0.992614727750844
    "F2" = 1.04824533977154
    "I2" = 1.03774609284961
    "J2" = 1.035860508230368
    "K2" = 1.04214509875397
    "L2" = 0.9955398523335997
    "M2" = 1.051005722764991
    "N2" = 1.015964129583316
    "B3" = 1.02
    "C3" = 1.03142819978128
    "D3" = 1.039898249494491
    "E3" = 0.9936372048519299
    "F3" = 1.048919317266124
    "I3" = 1.03787805516025
    "J3" = 1.036210691417483
    "K3" = 1.042494579353781
    "L3" = 0.9963617723202687
    "M3" = 1.051492064295724
    "N3" = 1.016080285508532
    "B4" = 1.02
    "C4" = 1.031886609894772
    "D4" = 1.040247363752442
    "E4" = 0.9942998659930998
    "F4" = 1.049356414942346
    "I4" = 1.037962630440587
    "J4" = 1.036437260402848
    "K4" = 1.042720634543939
    "L4" = 0.9968940712668347
    "M4" = 1.051807065354831
    "N4" = 1.016155426736826
    "B5" = 1.02
    "C5" = 1.032079488876482
    "D5" = 1.04039427994279
    "E5" = 0.994578699834602
    "F5" = 1.049540405323123
    "I5" = 1.037997990435024
    "J5" = 1.036532503177183
    "K5" = 1.042815647075079
    "L5" = 0.9971179600053012
    "M5" = 1.051939563151596
    "N5" = 1.016187010997427
    "B6" = 1.02
    "C6" = 1.032111883625499
    "D6" = 1.040418956470422
    "E6" = 0.994625531979634
    "F6" = 1.049571311803168
    "I6" = 1.038003916049491
    "J6" = 1.03654849441897
    "K6" = 1.04283159882051
    "L6" = 0.9971555583673455
    "M6" = 1.051961814246322
    "N6" = 1.016192313820554
    "B7" = 1.02
    "C7" = 1.031889186513858
    "D7" = 1.040249326272581
    "E7" = 0.994303590798249
    "F7" = 1.049358872512212
    "I7" = 1.037963103691537
    "J7" = 1.036438533069961
    "K7" = 1.042721904190596
    "L7" = 0.9968970624462089
    "M7" = 1.051808835518782
    "N7" = 1.016155848787929
    "B8" = 1.02
    "C8" = 1.030959741066795
    "D8" = 1.039541565921458
    "E8" = 0.9929600610674297
    "F8" = 1.048472907430698
    "I8" = 1.037790858030078
    "J8" = 1.035978858230029
    "K8" = 1.042263223517031
    "L8" = 0.9958175282591056
    "M8" = 1.051170019809789
    "N8" = 1.016003388766079
    "B9" = 1.02
    "C9" = 1.02932732533578
    "D9" = 1.038299307546956
    "E9" = 0.9906006454969559
    "F9" = 1.04691940091483
    "I9" = 1.037481145966678
    "J9" = 1.035168735545889
    "K9" = 1.041454402886394
    "L9" = 0.9939188001724441
    "M9" = 1.0500467648673
    "N9" = 1.015734606394287
    "B10" = 1.02
    "C10" = 1.028242794651588
    "D10" = 1.037474554097917
    "E10" = 0.989033133672735
    "F10" = 1.045889028829572
    "I10" = 1.037270554213103
    "J10" = 1.034628654926929
    "K10" = 1.040914892523868
    "L10" = 0.9926553831429383
    "M10" = 1.049299663199305
    "N10" = 1.015555359502377
    "B11" = 1.02
    "C11" = 1.027774096100146
    "D11" = 1.037118263240855
    "E11" = 0.988355674866747
    "F11" = 1.045444152355363
    "I11" = 1.037178399107527
    "J11" = 1.034394811266639
    "K11" = 1.04068122615278
    "L11" = 0.9921088820399291
    "M11" = 1.048976594370213
    "N11" = 1.015477735576271
    "B12" = 1.02
    "C12" = 1.02760013953226
    "D12" = 1.036986048170564
    "E12" = 0.9881042295826724
    "F12" = 1.045279100336422
    "I12" = 1.037144024042921
    "J12" = 1.03430795492128
    "K12" = 1.040594425279638
    "L12" = 0.9919059725120875
    "M12" = 1.048856658865047
    "N12" = 1.015448901742868
    "B13" = 1.02
    "C13" = 1.02763744748434
    "D13" = 1.037014402961507
    "E13" = 0.9881581567098651
    "F13" = 1.045314495696281
    "I13" = 1.037151404139565
    "J13" = 1.034326585717122
    "K13" = 1.040613044644962
    "L13" = 0.9919494934313052
    "M13" = 1.048882382392312
    "N13" = 1.015455086728876
    "B14" = 1.02
    "C14" = 1.027759713949056
    "D14" = 1.0371073316946
    "E14" = 0.9883348863814464
    "F14" = 1.045430505114045
    "I14" = 1.037175560598223
    "J14" = 1.034387631612603
    "K14" = 1.040674051295395
    "L14" = 0.9920921077337197
    "M14" = 1.048966679093644
    "N14" = 1.015475352177213
    "B15" = 1.02
    "C15" = 1.027835064829357
    "D15" = 1.037164605074252
    "E15" = 0.9884438009545853
    "F15" = 1.045502008244217
    "I15" = 1.037190425058605
    "J15" = 1.034425244522017
    "K15" = 1.040711638651488
    "L15" = 0.9921799884222134
    "M15" = 1.049018625960177
    "N15" = 1.015487838290256
    "B16" = 1.02
    "C16" = 1.028273919996713
    "D16" = 1.037498217674372
    "E16" = 0.9890781214508737
    "F16" = 1.045918581018568
    "I16" = 1.037276649928288
    "J16" = 1.034644174790141
    "K16" = 1.040930399160768
    "L16" = 0.9926916645766087
    "M16" = 1.049321113456569
    "N16" = 1.015560511001208
    "B17" = 1.02
    "C17" = 1.028549447516371
    "D17" = 1.037707708496989
    "E17" = 0.989476357848556
    "F17" = 1.046180230880452
    "I17" = 1.037330478051004
    "J17" = 1.034781508901242
    "K17" = 1.041067608200577
    "L17" = 0.9930127773699352
    "M17" = 1.049510972617415
    "N17" = 1.015606094643305
    "B18" = 1.02
    "C18" = 1.028710245751346
    "D18" = 1.03782998116526
    "E18" = 0.9897087662937556
    "F18" = 1.046332970298309
    "I18" = 1.037361781682258
    "J18" = 1.034861614901639
    "K18" = 1.041147634524765
    "L18" = 0.9932001317071769
    "M18" = 1.049621755725616
    "N18" = 1.015632681925752
    "B19" = 1.02
    "C19" = 1.028765088602616
    "D19" = 1.037871686527501
    "E19" = 0.9897880325774034
    "F19" = 1.04638507136425
    "I19" = 1.037372439544662
    "J19" = 1.034888929162499
    "K19" = 1.041174920475076
    "L19" = 0.9932640239640975
    "M19" = 1.049659536879986
    "N19" = 1.015641747337393
    "B20" = 1.02
    "C20" = 1.028519876934828
    "D20" = 1.03768522380213
    "E20" = 0.9894336180360679
    "F20" = 1.046152145537137
    "I20" = 1.037324712458119
    "J20" = 1.034766774102633
    "K20" = 1.04105288751587
    "L20" = 0.9929783193494215
    "M20" = 1.04949059822315
    "N20" = 1.015601204037985
    "B21" = 1.02
    "C21" = 1.027723705669707
    "D21" = 1.037079962970263
    "E21" = 0.9882828385668249
    "F21" = 1.045396337832288
    "I21" = 1.037168451107246
    "J21" = 1.034369655007437
    "K21" = 1.04065608653283
    "L21" = 0.9920501090198102
    "M21" = 1.048941853961471
    "N21" = 1.015469384527449
    "B22" = 1.02
    "C22" = 1.027223926148283
    "D22" = 1.036700148538545
    "E22" = 0.9875604150241495
    "F22" = 1.044922260408602
    "I22" = 1.037069367438044
    "J22" = 1.034119991990466
    "K22" = 1.040406563504805
    "L22" = 0.9914670000341481
    "M22" = 1.048597223161702
    "N22" = 1.015386499680789
    "B23" = 1.02
    "C23" = 1.027488791705152
    "D23" = 1.036901424790904
    "E23" = 0.9879432794643023
    "F23" = 1.04517347001941
    "I23" = 1.037121972537893
    "J23" = 1.034252340580559
    "K23" = 1.040538843546433
    "L23" = 0.991776070289318
    "M23" = 1.048779881209944
    "N23" = 1.015430438791648
    "B24" = 1.02
    "C24" = 1.028533238337114
    "D24" = 1.037695383420334
    "E24" = 0.9894529299347244
    "F24" = 1.046164835713662
    "I24" = 1.037327317970075
    "J24" = 1.034773432119932
    "K24" = 1.041059539176493
    "L24" = 0.9929938892766442
    "M24" = 1.049499804407031
    "N24" = 1.015603413894815
    "B25" = 1.02
    "C25" = 1.029748692762184
    "D25" = 1.038619867324398
    "E25" = 0.9912096547607049
    "F25" = 1.047320095646128
    "I25" = 1.037561943054127
    "J25" = 1.035378177305227
    "K25" = 1.041454402886394
    "L25" = 0.9944092447426414
    "M25" = 1.050336855285723
    "N25" = 1.01580410543714
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
